# Update view-count figures (column F) across the four worksheets to match
# the latest data pull (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 2275
$ws1.Range("F3").Value  = 345
$ws1.Range("F4").Value  = 181
$ws1.Range("F5").Value  = 185
$ws1.Range("F6").Value  = 346
$ws1.Range("F8").Value  = 709
$ws1.Range("F10").Value = 672
$ws1.Range("F14").Value = 973
$ws1.Range("F15").Value = 4807
$ws1.Range("F16").Value = 164
$ws1.Range("F18").Value = 26
$ws1.Range("F23").Value = 91
$ws1.Range("F24").Value = 14
$ws1.Range("F25").Value = 262

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 2777
$ws2.Range("F16").Value = 2511

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 403

# --- Sheet "全部类型" (All Types, merged view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 2275
$ws4.Range("F7").Value  = 403
$ws4.Range("F8").Value  = 345
$ws4.Range("F9").Value  = 181
$ws4.Range("F10").Value = 185
$ws4.Range("F11").Value = 346
$ws4.Range("F17").Value = 709
$ws4.Range("F19").Value = 672
$ws4.Range("F23").Value = 973
$ws4.Range("F24").Value = 4807
$ws4.Range("F26").Value = 2777
$ws4.Range("F30").Value = 164
$ws4.Range("F32").Value = 26
$ws4.Range("F41").Value = 91
$ws4.Range("F42").Value = 14
$ws4.Range("F43").Value = 262
$ws4.Range("F45").Value = 2511
